$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: round coordinates, drop start/end time values ---
$ws.Range("Q2").Value = 492738
$ws.Range("R2").Value = 6845111
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# --- Row 3: round coordinates, drop start/end time values ---
$ws.Range("Q3").Value = 492845
$ws.Range("R3").Value = 6845301
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# --- Rows 4 & 5: the two species records are swapped (row 4 <-> row 5),
#     their coordinates rounded, and start/end time values dropped ---
$ws.Range("A4").Value = 112181727
$ws.Range("B4").Value = 78578
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6458
$ws.Range("F4").Value = "Lunglav"
$ws.Range("G4").Value = "Lobaria pulmonaria"
$ws.Range("H4").Value = "(L.) Hoffm."
$ws.Range("Q4").Value = 493005
$ws.Range("R4").Value = 6845384
$ws.Range("Y4").Value = "'2023-07-04"
$ws.Range("Y4").Style = "Normal"
$ws.Range("Z4").ClearContents()
$ws.Range("AA4").Value = "'2023-07-04"
$ws.Range("AA4").Style = "Normal"
$ws.Range("AB4").ClearContents()

$ws.Range("A5").Value = 112182724
$ws.Range("B5").Value = 95538
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 221941
$ws.Range("F5").Value = "Plattlummer"
$ws.Range("G5").Value = "Lycopodium complanatum"
$ws.Range("H5").Value = "L."
$ws.Range("Q5").Value = 493165
$ws.Range("R5").Value = 6845494
$ws.Range("Y5").Value = "'2023-07-06"
$ws.Range("Y5").Style = "Normal"
$ws.Range("Z5").ClearContents()
$ws.Range("AA5").Value = "'2023-07-06"
$ws.Range("AA5").Style = "Normal"
$ws.Range("AB5").ClearContents()

# --- Row 6: round coordinates, drop start/end time values ---
$ws.Range("Q6").Value = 492536
$ws.Range("R6").Value = 6845328
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
